$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 201; $r -le 312; $r++) {
    $ws.Cells.Item($r, 2).Value = "yes"
}
